# Update the three numeric values in column C on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.3
$ws.Range("C3").Value = 2.3
$ws.Range("C4").Value = 3.3

# Move/collapse the active selection to F8 (single cell)
$ws.Range("F8").Select()
